# Update NATMI TPM-derived values on the active sheet.
# Only the Receptor average/total expression (M, N) are "new TPM" inputs;
# the derived specificity / edge-weight columns (O-T) are recomputed from
# them (and are rewritten here with their already-computed target values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.6436386666666666
$ws.Range("N2").Value = 1.930916
$ws.Range("O2").Value = 0.3148741541159968
$ws.Range("P2").Value = 0.3148741541159968
$ws.Range("Q2").Value = 0.1355022448462222
$ws.Range("R2").Value = 1.219520203616
$ws.Range("S2").Value = 0.3148741541159968
$ws.Range("T2").Value = 0.3148741541159968

# Row 3 (receptor expression itself unchanged; only specificities shift
# because the per-cluster sums moved)
$ws.Range("O3").Value = 0.2664627967585631
$ws.Range("P3").Value = 0.266462796758563
$ws.Range("S3").Value = 0.2664627967585631
$ws.Range("T3").Value = 0.266462796758563

# Row 4
$ws.Range("M4").Value = 0.1891833333333333
$ws.Range("N4").Value = 0.56755
$ws.Range("O4").Value = 0.09255028502976516
$ws.Range("P4").Value = 0.09255028502976514
$ws.Range("Q4").Value = 0.03982788431111112
$ws.Range("R4").Value = 0.3584509588
$ws.Range("S4").Value = 0.09255028502976516
$ws.Range("T4").Value = 0.09255028502976514

# Row 5
$ws.Range("M5").Value = 0.3089366666666667
$ws.Range("N5").Value = 0.92681
$ws.Range("O5").Value = 0.1511347540629665
$ws.Range("P5").Value = 0.1511347540629665
$ws.Range("Q5").Value = 0.0650389947288889
$ws.Range("R5").Value = 0.58535095256
$ws.Range("S5").Value = 0.1511347540629665
$ws.Range("T5").Value = 0.1511347540629665

# Row 6
$ws.Range("M6").Value = 0.1497823333333334
$ws.Range("N6").Value = 0.4493470000000001
$ws.Range("O6").Value = 0.07327494128670582
$ws.Range("P6").Value = 0.07327494128670581
$ws.Range("Q6").Value = 0.03153297565244445
$ws.Range("R6").Value = 0.2837967808720001
$ws.Range("S6").Value = 0.07327494128670582
$ws.Range("T6").Value = 0.07327494128670581

# Row 7
$ws.Range("M7").Value = 0.2078926666666666
$ws.Range("N7").Value = 0.623678
$ws.Range("O7").Value = 0.1017030687460028
$ws.Range("P7").Value = 0.1017030687460028
$ws.Range("Q7").Value = 0.04376667294755555
$ws.Range("R7").Value = 0.393900056528
$ws.Range("S7").Value = 0.1017030687460028
$ws.Range("T7").Value = 0.1017030687460028
